# "Simple storage to EEPROM"
#
# The "Trigs off" row (row 25 on sheet "1k") is removed — the design now
# stores a single combined "Trigs" vector instead of separate "notes on" /
# "notes off" vectors. Removing that row shifts every row below it up by
# one, and the row that used to hold "Trigs on" (row 24) is relabeled to
# just "Trigs".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1k")

# The TOTAL formula (old row 31, new row 30) referenced E25 ("Trigs off"
# byte count) directly. Strip that term out before the row disappears so
# we don't end up with a #REF! error; the remaining terms simply shift
# their cell references up automatically when the row is deleted.
$ws.Range("E31").Formula = "=E11+B29*(E17+E23+E24)+E28"

# Delete the whole "Trigs off" row (row 25). Everything below (rows 26-32)
# shifts up to become rows 25-31.
$ws.Rows("25:25").Delete()

# The conditional formatting rule ("highlight if < 0") was anchored to the
# old leftover-bytes cell E32; move it to its new location E31 (this keeps
# reusing the existing dxf style instead of creating a new one).
$fc = $ws.Range("E32").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("E31"))

# Rename what used to be "Trigs on" (now row 24) to simply "Trigs".
$ws.Range("A24").Value = "Trigs"

# Match the author's final selection/cursor position.
$ws.Range("O25").Select()
